# Planeacion cronograma proyecto.xlsx - status column updates + view state
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Example - Project Plan Template")

# --- Update task status values in column D (sheet1 / "Example - Project Plan Template") ---
# 3.x tasks: In Progress -> Complete
$ws.Range("D31").Value = "Complete"
$ws.Range("D32").Value = "Complete"
$ws.Range("D33").Value = "Complete"

# 4.1 task: In Progress -> Complete
$ws.Range("D35").Value = "Complete"

# 4.3, 5.1, 5.2 tasks: Not Started -> Complete
$ws.Range("D37").Value = "Complete"
$ws.Range("D39").Value = "Complete"
$ws.Range("D40").Value = "Complete"

# 6.1, 6.2, 6.3 tasks: Not Started -> In Progress
$ws.Range("D42").Value = "In Progress"
$ws.Range("D43").Value = "In Progress"
$ws.Range("D44").Value = "In Progress"

# --- Update sheet view state (zoom + frozen-pane selection) ---
$win = $excel.ActiveWindow

# Re-freeze the header row (ySplit=1) and zoom to 85%, matching the saved view.
$win.FreezePanes = $false
$ws.Range("A2").Select() | Out-Null
$win.FreezePanes = $true

$win.Zoom = 85

# Move the active selection to the cell that was selected when the file was saved.
$ws.Range("D44").Select() | Out-Null
